# Add season-record columns (Wins, Losses, Ties) to the player table.
# Mirrors the author's fix: the scraper previously only pulled team
# statistics, not the season win/loss/tie record, so these three new
# columns carry that record for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the
# three new header cells so they pick up the same bold/bordered/centered
# style used by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player (rows 2-46) shares the same team season record for 2023:
# 90 wins, 72 losses, 0 ties.
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 90
    $ws.Cells.Item($r, 31).Value = 72
    $ws.Cells.Item($r, 32).Value = 0
}
